# Apply the data updates described in the commit:
# "Made changes in the data file"
#
# Row 3: C3 "gkjgfkj" -> "gkjgfkj123", D3 "kjdfkjdf" -> "kjdfkjdf123"
# Row 4: A4 "asdfg" -> "Dokuparthi", B4 "asdfg@asdf.com" -> "dokuparthi@gmail.com",
#        C4 "gfkjkjj" -> "123Dokuparthi", D4 "ljnkdv" -> "123fdsdfdf"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for how new shared strings get appended to xl/sharedStrings.xml,
# so apply the edits in the same order the original author made them:
# A4, then C3, D3, then B4, C4, D4.
$ws.Range("A4").Value = "Dokuparthi"
$ws.Range("C3").Value = "gkjgfkj123"
$ws.Range("D3").Value = "kjdfkjdf123"
$ws.Range("B4").Value = "dokuparthi@gmail.com"
$ws.Range("C4").Value = "123Dokuparthi"
$ws.Range("D4").Value = "123fdsdfdf"

# Update the active selection to match the saved view state (C3)
$ws.Range("C3").Select()
